# Add team record (Wins / Losses / Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 53            # last data row (row 1 is the header)
$firstDataRow = 2

# --- Header row (row 1): new columns AD, AE, AF ------------------------------
# Copy the formatting of the existing last header cell (AC1) onto the three
# new header cells so they keep the bold/centered/bordered look of the rest
# of the header row, then set their text.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# --- Data rows (2-53): Wins = 65, Losses = 97, Ties = 0 ----------------------
$rowCount = $lastRow - $firstDataRow + 1

$wins = New-Object 'object[,]' $rowCount,1
$losses = New-Object 'object[,]' $rowCount,1
$ties = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $wins[$i,0] = 65
    $losses[$i,0] = 97
    $ties[$i,0] = 0
}

$ws.Range("AD$firstDataRow`:AD$lastRow").Value2 = $wins
$ws.Range("AE$firstDataRow`:AE$lastRow").Value2 = $losses
$ws.Range("AF$firstDataRow`:AF$lastRow").Value2 = $ties
